# Prep new data 2024: append 8 new GHG indicator rows to the "Series" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append below the existing data (rows 2..30), starting at row 31.
# Columns: A=Code, B=Name, C=Source Note, D=Derivation Rule, E=Weights Indicator,
#          F=Aggregation Rule, G=Indicator Type
$newRows = @(
    @{ Code = "EN.GHG.ALL.PC.CE.AR5";      Name = "Total greenhouse gas emissions per capita excluding LULUCF (t CO2e/capita)" },
    @{ Code = "EN.GHG.CO2.PC.CE.AR5";      Name = "Carbon dioxide (CO2) emissions excluding LULUCF per capita (t CO2e/capita)" },
    @{ Code = "EN.GHG.TOT.ZG.AR5";         Name = "Total greenhouse gas emissions excluding LULUCF (% change from 1990)" },
    @{ Code = "EN.GHG.CO2.ZG.AR5";         Name = "Carbon dioxide (CO2) emissions (total) excluding LULUCF (% change from 1990)" },
    @{ Code = "EN.GHG.CH4.ZG.AR5";         Name = "Methane (CH4) emissions (total) excluding LULUCF (% change from 1990)" },
    @{ Code = "EN.GHG.N2O.ZG.AR5";         Name = "Nitrous oxide (N2O) emissions (total) excluding LULUCF (% change from 1990)" },
    @{ Code = "EN.GHG.CO2.RT.GDP.KD";      Name = "Carbon intensity of GDP (kg CO2e per 2021 US$ of GDP)" },
    @{ Code = "EN.GHG.CO2.RT.GDP.PP.KD";   Name = "Carbon intensity of GDP (kg CO2e per 2021 PPP $)" }
)

$sourceNote = "EDGAR - Emissions Database for Global Atmospheric Research"
$indicatorType = "Basic"

$startRow = 31

# First pass: fill in Code (A) and Name (B) for every new row - this mirrors
# how the shared-string table was actually built in the source edit.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row.Code
    $ws.Cells.Item($r, 2).Value = $row.Name
}

# Second pass: fill in Source Note (C) and Indicator Type (G) for every row.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 3).Value = $sourceNote
    $ws.Cells.Item($r, 7).Value = $indicatorType
}

# Update the view so the newly added rows are visible / selected, matching
# the author's saved window state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I39").Select()
